# Generate Report for Handoff
# Updates the "d4543783-963b-44e2-9069-77ac2bffa0e0.md" row (row 3) on every
# sheet: status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", the handoff timestamps advance, and the zh-cn / de-de
# sheets gain an "Error Detail" message describing why the handback file is
# stale. The Error Detail column is also widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d16c5df22cc9a219a7af71669bf6a93cefc8368f/e2e/d4543783-963b-44e2-9069-77ac2bffa0e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de828b4372e52008b150b328106b642c97862f33/e2e/d4543783-963b-44e2-9069-77ac2bffa0e0.md."

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $readyForHandoff          # zh-cn status
$overview.Range("F3").Value = $readyForHandoff          # de-de status
$overview.Range("G3").Value = "2016-08-27 00:48:35"     # Latest HO Xliff Generate Date

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $readyForHandoff               # Status
$zhcn.Range("H3").Value = "2016-08-27 00:48:31"          # Latest Handoff Datetime
$zhcn.Range("P3").Value = $errorDetail                   # Error Detail
$zhcn.Range("P1").ColumnWidth = 235/6                    # widen Error Detail column to 40

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyForHandoff               # Status
$dede.Range("H3").Value = "2016-08-27 00:48:35"          # Latest Handoff Datetime
$dede.Range("P3").Value = $errorDetail                   # Error Detail
$dede.Range("P1").ColumnWidth = 235/6                    # widen Error Detail column to 40
